$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 44737
$ws.Range("A11").NumberFormat = "m/d/yyyy"
$ws.Range("B11").Value = "Octavio Lucardi Fierro"
$ws.Range("C11").Value = "CHICOS AYUDAAAA QUE ALGUIEN MAS HAGA SPRITES PORQUE ES HORRIBLE "
$ws.Range("D11").Value = "juan tien 2 tipos de ataque, un apuñalamiento a poca distancia y otro como un barrido, al cual no le funciona la animacion porque UNITY no me deja hacer la cosas igual que en los tutoriales"

$ws.Range("A11:D11").RowHeight = 60
$ws.Range("C10").Select()
$ws.Range("C11").Select()
